$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.071.68"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.017.90"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "227.36"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'55.80"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("E11").Value = "  -3.51%  "
$ws.Range("D12").Value = "2.317.42"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").Value = "14.18"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "5.16"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "2.027.44"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "37.026.60"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "68.87"
$ws.Range("D21").Value = "0.0₃0815"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "222.89"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  +3.00%  "
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("D26").Value = "163.38"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("E27").Value = "  -4.40%  "
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "4.45"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").Value = "2.34"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "3.16"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "5.45"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("D40").Value = "1.469.43"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").Value = "0.0212"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").Value = "94.38"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").Value = "'0.0910"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").Value = "'16.20"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "4.15"
$ws.Range("E46").Value = "  +13.19%  "
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").Value = "7.04"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").Value = "2.205.44"
$ws.Range("E51").Value = "  -0.82%  "
